$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=45693.82777777778; B=26.88},
    @{Row=3;  A=45695.51388888889; B=26.21},
    @{Row=4;  A=45700.33194444444; B=27.05},
    @{Row=5;  A=45704.49930555555; B=27.91},
    @{Row=6;  A=45709.9125;        B=25.41},
    @{Row=7;  A=45710.32013888889; B=27.99},
    @{Row=8;  A=45711.50347222222; B=26.9},
    @{Row=9;  A=45711.75069444445; B=25.43},
    @{Row=10; A=45715.97152777778; B=24.91},
    @{Row=11; A=45716.13888888889; B=30.31}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
}
